# Auto-generated script to apply cryptos.xlsx cell-value updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.585.31"
$ws.Range("E2").Value = "  +0.19%  "
$ws.Range("D3").Value = "1.647.82"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.61"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.63%  "
$ws.Range("E6").Value = "  +4.94%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("E8").Value = "  -1.21%  "
$ws.Range("E9").Value = "  -1.57%  "
$ws.Range("E10").Value = "  -1.12%  "
$ws.Range("E11").Value = "  +1.19%  "
$ws.Range("D12").Value = "1.880.46"
$ws.Range("E12").Value = "  -0.43%  "
$ws.Range("D13").Value = "1.647.29"
$ws.Range("E13").Value = "  -0.43%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.583"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.61%  "
$ws.Range("E15").Value = "  -2.32%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.50"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.92%  "
$ws.Range("D17").Value = "27.550.27"
$ws.Range("E17").Value = "  +0.11%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "231.58"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.65%  "
$ws.Range("E19").Value = "  -0.73%  "
$ws.Range("E20").Value = "  +0.37%  "
$ws.Range("E21").Value = "  -0.06%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.33"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.13%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.78"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.38%  "
$ws.Range("E24").Value = "  -1.60%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "149.07"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.25%  "
$ws.Range("E26").Value = "  -2.55%  "
$ws.Range("E27").Value = "  +1.53%  "
$ws.Range("E28").Value = "  -0.07%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.61"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.85%  "
$ws.Range("E30").Value = "  -2.03%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0486"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.99%  "
$ws.Range("E32").Value = "  -0.58%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.17"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Value = "1.422.51"
$ws.Range("E34").Value = "  -1.74%  "
$ws.Range("E35").Value = "  +2.67%  "
$ws.Range("E36").Value = "  -0.27%  "
$ws.Range("E37").Value = "  -0.09%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.884"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.86%  "
$ws.Range("E39").Value = "  -2.27%  "
$ws.Range("E40").Value = "  +0.17%  "
$ws.Range("E41").Value = "  -0.06%  "
$ws.Range("E42").Value = "  +2.79%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.54"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.33%  "
$ws.Range("E44").Value = "  +1.43%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "65.15"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.81%  "
$ws.Range("D46").Value = "1.790.16"
$ws.Range("E46").Value = "  -0.33%  "
$ws.Range("E47").Value = "  -1.39%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "88.23"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.22%  "
$ws.Range("D49").Value = "0.0₆0107"
$ws.Range("E49").Value = "  +1.03%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0996"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.17%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.80"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.13%  "
